$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header "muc_do_cau_hoi" in G6 (keeps existing cell style) ---
$ws.Range("G6").Value = "muc_do_cau_hoi"

# --- New note row (A4): explains the muc_do_cau_hoi column ---
$noteText = "Cột muc_do_cau_hoi yêu cầu phải nhập đúng như trong hệ thống đã ghi để tránh bị lỗi khi nhập vào"
$ws.Range("A4").Value = $noteText

# "Cột " (chars 1-4) stays default formatting.
# "muc_do_cau_hoi" (chars 5-18) is bold, Arial 10, black.
$ws.Range("A4").Characters(5, 14).Font.Bold = $true
$ws.Range("A4").Characters(5, 14).Font.Size = 10
$ws.Range("A4").Characters(5, 14).Font.Color = 0
$ws.Range("A4").Characters(5, 14).Font.Name = "Arial"

# Remaining text (chars 19-96) is regular, Arial 10, black.
$ws.Range("A4").Characters(19, 78).Font.Size = 10
$ws.Range("A4").Characters(19, 78).Font.Color = 0
$ws.Range("A4").Characters(19, 78).Font.Name = "Arial"

# --- Column width adjustments (C, D, E, F, G got narrower) ---
$ws.Columns("C").ColumnWidth = 31.833333333333332
$ws.Columns("D").ColumnWidth = 27.166666666666668
$ws.Columns("E").ColumnWidth = 26.166666666666668
$ws.Columns("F").ColumnWidth = 33.0
$ws.Columns("G").ColumnWidth = 15.666666666666666

# --- Move the active selection from A3:N3 to B9 ---
$ws.Range("B9").Select()
